# Updated capital structure database
# - Removes the "Prime Media Holdings, Inc. (PSE:PRIM)" row (row 5)
# - Refreshes computed metrics for the remaining three rows
# - Row 3 / Row 4 company identities are swapped (PBB now row 3, AUB now row 4)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the last data row (Prime Media Holdings, Inc. (PSE:PRIM)); this also
# shrinks the sheet dimension from A1:AQ5 down to A1:AQ4.
$ws.Rows(5).Delete()

# --- row2 ---
$ws.Range("A2").Value = "Philippines"
$ws.Range("B2").Value = "'2"
$ws.Range("C2").Value = "Banks (Regional)"
$ws.Range("D2").Value = 0.164
$ws.Range("E2").Value = 0.258
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 122.5
$ws.Range("L2").Value = 0.361890694239291
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 0.03199488081906895
$ws.Range("O2").Value = 0.163265306122449
$ws.Range("P2").Value = 20
$ws.Range("Q2").Value = 0.03199488081906895
$ws.Range("R2").Value = 0.163265306122449
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 286.4
$ws.Range("V2").Value = 0.4581666933290673
$ws.Range("W2").Value = 0.1399288701278692
$ws.Range("X2").Value = 0.05864253609723899
$ws.Range("Y2").Value = 0.08128633403063021
$ws.Range("Z2").Value = 0.3781066741133762
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.04668753046324002
$ws.Range("AC2").Value = -0.04668753046324002
$ws.Range("AD2").Value = 390.2
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 390.2
$ws.Range("AG2").Value = 103.8
$ws.Range("AH2").Value = 0.384319905446666
$ws.Range("AI2").Value = 0.2746920098556846
$ws.Range("AJ2").Value = 0.1424063657566196
$ws.Range("AK2").Value = 0.09152632043029715
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0

# --- row3 ---
$ws.Range("A3").Value = "Philippines"
$ws.Range("B3").Value = "Philippine Business Bank, Inc. (PSE:PBB)"
$ws.Range("C3").Value = "Banks (Regional)"
$ws.Range("D3").Value = 0.167
$ws.Range("E3").Value = 0.237
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 31.5
$ws.Range("L3").Value = 0.2858439201451906
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("S3").Value = 0
$ws.Range("U3").Value = 213.1
$ws.Range("V3").Value = 1.242565597667638
$ws.Range("W3").Value = 0.1353674258702192
$ws.Range("X3").Value = 0.05478799770869061
$ws.Range("Y3").Value = 0.08057942816152855
$ws.Range("Z3").Value = 0.3946284691136975
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.04634749180802263
$ws.Range("AC3").Value = -0.04634749180802263
$ws.Range("AD3").Value = 68.5
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 68.5
$ws.Range("AG3").Value = -144.6
$ws.Range("AH3").Value = 0.2854166666666667
$ws.Range("AI3").Value = 0.1894358407079646
$ws.Range("AJ3").Value = -5.375464684014869
$ws.Range("AK3").Value = -0.9737373737373735
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = 0

# --- row4 ---
$ws.Range("A4").Value = "Philippines"
$ws.Range("B4").Value = "Asia United Bank Corporation (PSE:AUB)"
$ws.Range("C4").Value = "Banks (Regional)"
$ws.Range("D4").Value = 0.161
$ws.Range("E4").Value = 0.279
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 91
$ws.Range("L4").Value = 0.398598335523434
$ws.Range("M4").Value = 20
$ws.Range("N4").Value = 0.04409171075837742
$ws.Range("O4").Value = 0.2197802197802198
$ws.Range("P4").Value = 20
$ws.Range("Q4").Value = 0.04409171075837742
$ws.Range("R4").Value = 0.2197802197802198
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 73.3
$ws.Range("V4").Value = 0.1615961199294532
$ws.Range("W4").Value = 0.1444903143855192
$ws.Range("X4").Value = 0.06249707448578737
$ws.Range("Y4").Value = 0.08199323989973187
$ws.Range("Z4").Value = 0.3706168831168831
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.04702756911845742
$ws.Range("AC4").Value = -0.04702756911845742
$ws.Range("AD4").Value = 321.7
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 321.7
$ws.Range("AG4").Value = 248.4
$ws.Range("AH4").Value = 0.4149361537469367
$ws.Range("AI4").Value = 0.3038058362451601
$ws.Range("AJ4").Value = 0.3538461538461538
$ws.Range("AK4").Value = 0.2520292207792207
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0

# Cells that no longer carry a value in the refreshed export (debt_ebitda /
# net_debt_ebitda were dropped for every remaining row).
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()

# buybacks_cash_returned is no longer reported for row 3 (Philippine Business
# Bank) in the refreshed export.
$ws.Range("T3").ClearContents()
